# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" sheets to reflect refreshed scrape output.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 72
$ws1.Range("F5").Value = 74
$ws1.Range("F9").Value = 1095
$ws1.Range("F10").Value = 15322
$ws1.Range("F11").Value = 213
$ws1.Range("F12").Value = 159
$ws1.Range("F14").Value = 6048
$ws1.Range("F17").Value = 62
$ws1.Range("F18").Value = 100
$ws1.Range("F24").Value = 845
$ws1.Range("F25").Value = 11
$ws1.Range("F26").Value = 4973
$ws1.Range("F28").Value = 10915
$ws1.Range("F30").Value = 5
$ws1.Range("F31").Value = 103
$ws1.Range("F32").Value = 147
$ws1.Range("F33").Value = 3780
$ws1.Range("F34").Value = 257
$ws1.Range("F35").Value = 73

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 72
$ws4.Range("F6").Value = 74
$ws4.Range("F10").Value = 1095
$ws4.Range("F11").Value = 15322
$ws4.Range("F12").Value = 213
$ws4.Range("F13").Value = 159
$ws4.Range("F15").Value = 6048
$ws4.Range("F18").Value = 62
$ws4.Range("F19").Value = 100
$ws4.Range("F25").Value = 845
$ws4.Range("F26").Value = 11
$ws4.Range("F27").Value = 4973
$ws4.Range("F30").Value = 10915
$ws4.Range("F32").Value = 5
$ws4.Range("F33").Value = 103
$ws4.Range("F34").Value = 147
$ws4.Range("F35").Value = 3780
$ws4.Range("F36").Value = 257
$ws4.Range("F37").Value = 73
